# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns for rows 2-51.
# Values that look like plain decimal numbers (e.g. "0.9997") are written with a
# leading apostrophe so Excel stores them as text (matching the source inlineStr
# cells) instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $v = $text
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Looks like a plain decimal number: prefix with an apostrophe so Excel
        # keeps it as text (matching the source inlineStr cell) instead of
        # silently re-typing it as a numeric value.
        $v = "'" + $text
    }
    $ws.Cells.Item($row, $col).Value = $v
}

Set-TextValue 2 4 '29.988.17'
Set-TextValue 2 5 '  -0.90%  '

Set-TextValue 3 4 '1.902.39'
Set-TextValue 3 5 '  -0.70%  '

Set-TextValue 4 5 '  -0.57%  '

Set-TextValue 5 4 '0.7452'
Set-TextValue 5 5 '  +1.06%  '

Set-TextValue 6 4 '241.24'
Set-TextValue 6 5 '  -0.41%  '

Set-TextValue 7 4 '0.9997'

Set-TextValue 8 4 '0.3069'
Set-TextValue 8 5 '  -1.62%  '

Set-TextValue 9 4 '25.52'
Set-TextValue 9 5 '  -5.73%  '

Set-TextValue 10 4 '0.06893'
Set-TextValue 10 5 '  -0.42%  '

Set-TextValue 11 4 '0.08012'
Set-TextValue 11 5 '  +0.20%  '

Set-TextValue 12 4 '0.7545'
Set-TextValue 12 5 '  -1.46%  '

Set-TextValue 13 4 '1.901.34'
Set-TextValue 13 5 '  -0.52%  '

Set-TextValue 14 4 '5.265'
Set-TextValue 14 5 '  -0.25%  '

Set-TextValue 15 4 '91.48'
Set-TextValue 15 5 '  +0.44%  '

Set-TextValue 16 4 '6.156'
Set-TextValue 16 5 '  +5.72%  '

Set-TextValue 17 4 '30.001.10'
Set-TextValue 17 5 '  -0.73%  '

Set-TextValue 18 4 '14.03'
Set-TextValue 18 5 '  -0.66%  '

Set-TextValue 19 4 '0.000007751'
Set-TextValue 19 5 '  -0.84%  '

Set-TextValue 20 4 '237.58'
Set-TextValue 20 5 '  -3.41%  '

Set-TextValue 21 4 '0.9997'

Set-TextValue 22 4 '2.152.64'
Set-TextValue 22 5 '  -0.06%  '

Set-TextValue 23 4 '0.9995'
Set-TextValue 23 5 '  -0.63%  '

Set-TextValue 24 4 '7.100'
Set-TextValue 24 5 '  +8.04%  '

Set-TextValue 25 4 '9.310'
Set-TextValue 25 5 '  -0.40%  '

Set-TextValue 26 4 '166.30'
Set-TextValue 26 5 '  +0.78%  '

Set-TextValue 27 4 '18.80'
Set-TextValue 27 5 '  -0.10%  '

Set-TextValue 28 5 '  -0.62%  '

Set-TextValue 29 4 '2.056'
Set-TextValue 29 5 '  -3.50%  '

Set-TextValue 30 4 '1.343'
Set-TextValue 30 5 '  -0.94%  '

Set-TextValue 31 5 '  -0.89%  '

Set-TextValue 32 4 '4.304'
Set-TextValue 32 5 '  -0.33%  '

Set-TextValue 33 4 '4.047'
Set-TextValue 33 5 '  +0.10%  '

Set-TextValue 34 4 '0.05423'
Set-TextValue 34 5 '  +5.41%  '

Set-TextValue 35 4 '1.284'
Set-TextValue 35 5 '  -0.58%  '

Set-TextValue 36 4 '0.7379'
Set-TextValue 36 5 '  -0.38%  '

Set-TextValue 37 4 '2.724'
Set-TextValue 37 5 '  -1.90%  '

Set-TextValue 38 4 '0.01943'
Set-TextValue 38 5 '  +0.79%  '

Set-TextValue 39 4 '2.770'
Set-TextValue 39 5 '  -0.24%  '

Set-TextValue 40 4 '6.235'
Set-TextValue 40 5 '  -2.56%  '

Set-TextValue 41 4 '0.4455'
Set-TextValue 41 5 '  +0.44%  '

Set-TextValue 42 4 '72.60'
Set-TextValue 42 5 '  -3.73%  '

Set-TextValue 43 4 '1.944'
Set-TextValue 43 5 '  +0.97%  '

Set-TextValue 44 4 '0.9997'
Set-TextValue 44 5 '  -0.54%  '

Set-TextValue 45 4 '0.8312'
Set-TextValue 45 5 '  -0.59%  '

Set-TextValue 46 4 '7.672'
Set-TextValue 46 5 '  +1.73%  '

Set-TextValue 47 4 '101.52'
Set-TextValue 47 5 '  +0.72%  '

Set-TextValue 48 4 '9.858'
Set-TextValue 48 5 '  +1.25%  '

Set-TextValue 49 4 '2.057.96'
Set-TextValue 49 5 '  -0.54%  '

Set-TextValue 50 4 '36.55'
Set-TextValue 50 5 '  -0.98%  '

Set-TextValue 51 4 '0.1163'
Set-TextValue 51 5 '  -3.09%  '
